# Auto-generated Excel COM-interop script
# Applies cached numeric-value updates (H:N columns) produced by the
# scheduled market-price refresh run, per the supplied OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Cells.Item(9, 8).Value = 249.25
$ws.Cells.Item(9, 10).Value = 238.5
$ws.Cells.Item(9, 12).Value = 238.5
$ws.Cells.Item(9, 14).Value = -576.5
# Row 12
$ws.Cells.Item(12, 8).Value = 103.333336
$ws.Cells.Item(12, 9).Value = 105
$ws.Cells.Item(12, 11).Value = 105
$ws.Cells.Item(12, 13).Value = 65
# Row 17
$ws.Cells.Item(17, 8).Value = 1565.1875
$ws.Cells.Item(17, 10).Value = 1598.6666
$ws.Cells.Item(17, 12).Value = 4795.9998
$ws.Cells.Item(17, 14).Value = -5131.9998
# Row 80
$ws.Cells.Item(80, 8).Value = 492.43478
$ws.Cells.Item(80, 10).Value = 632.44446
$ws.Cells.Item(80, 12).Value = 1897.33338
$ws.Cells.Item(80, 14).Value = -3893.33338
# Row 83
$ws.Cells.Item(83, 8).Value = 492.43478
$ws.Cells.Item(83, 10).Value = 632.44446
$ws.Cells.Item(83, 12).Value = 5692.00014
$ws.Cells.Item(83, 14).Value = -15676.00014
# Row 98
$ws.Cells.Item(98, 8).Value = 2003.5714
$ws.Cells.Item(98, 9).Value = 2576.6667
$ws.Cells.Item(98, 11).Value = 2576.6667
$ws.Cells.Item(98, 13).Value = -1078.6667
# Row 122
$ws.Cells.Item(122, 8).Value = 2003.5714
$ws.Cells.Item(122, 9).Value = 2576.6667
$ws.Cells.Item(122, 11).Value = 7730.000100000001
$ws.Cells.Item(122, 13).Value = -5280.000100000001
# Row 132
$ws.Cells.Item(132, 8).Value = 2048.5833
$ws.Cells.Item(132, 9).Value = 2123.7778
$ws.Cells.Item(132, 11).Value = 6371.3334
$ws.Cells.Item(132, 13).Value = -3841.3334
# Row 135
$ws.Cells.Item(135, 8).Value = 1285.8334
$ws.Cells.Item(135, 9).Value = 893
$ws.Cells.Item(135, 11).Value = 8037
$ws.Cells.Item(135, 13).Value = -5502
# Row 137
$ws.Cells.Item(137, 8).Value = 1300.625
$ws.Cells.Item(137, 10).Value = 1980
$ws.Cells.Item(137, 12).Value = 5940
$ws.Cells.Item(137, 14).Value = -11040

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Cells.Item(24, 8).Value = 27495
$ws.Cells.Item(24, 10).Value = 27495
$ws.Cells.Item(24, 12).Value = 27495
$ws.Cells.Item(24, 14).Value = -28243
# Row 61
$ws.Cells.Item(61, 8).Value = 2643.2222
$ws.Cells.Item(61, 9).Value = 2643.2222
$ws.Cells.Item(61, 11).Value = 2643.2222
$ws.Cells.Item(61, 13).Value = -2431.2222
# Row 100
$ws.Cells.Item(100, 8).Value = 27495
$ws.Cells.Item(100, 10).Value = 27495
$ws.Cells.Item(100, 12).Value = 27495
$ws.Cells.Item(100, 14).Value = -29659
# Row 132
$ws.Cells.Item(132, 8).Value = 1285.4166
$ws.Cells.Item(132, 9).Value = 1215.2222
$ws.Cells.Item(132, 11).Value = 3645.6666
$ws.Cells.Item(132, 13).Value = -1115.6666
# Row 136
$ws.Cells.Item(136, 8).Value = 2643.2222
$ws.Cells.Item(136, 9).Value = 2643.2222
$ws.Cells.Item(136, 11).Value = 7929.6666
$ws.Cells.Item(136, 13).Value = -5379.6666

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 2997
$ws.Cells.Item(94, 9).Value = 2996
$ws.Cells.Item(94, 11).Value = 2996
$ws.Cells.Item(94, 13).Value = -2545
# Row 105
$ws.Cells.Item(105, 8).Value = 1936.9565
$ws.Cells.Item(105, 9).Value = 1692.9048
$ws.Cells.Item(105, 10).Value = 4499.5
$ws.Cells.Item(105, 11).Value = 1692.9048
$ws.Cells.Item(105, 12).Value = 4499.5
$ws.Cells.Item(105, 13).Value = 54.09519999999998
$ws.Cells.Item(105, 14).Value = -7993.5
# Row 134
$ws.Cells.Item(134, 8).Value = 1711.2927
$ws.Cells.Item(134, 9).Value = 983.8889
$ws.Cells.Item(134, 10).Value = 3114.1428
$ws.Cells.Item(134, 11).Value = 2951.6667
$ws.Cells.Item(134, 12).Value = 9342.428400000001
$ws.Cells.Item(134, 13).Value = -416.6667000000002
$ws.Cells.Item(134, 14).Value = -14412.4284

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2805.56
$ws.Cells.Item(31, 9).Value = 1939.4736
$ws.Cells.Item(31, 10).Value = 5548.1665
$ws.Cells.Item(31, 11).Value = 1939.4736
$ws.Cells.Item(31, 12).Value = 5548.1665
$ws.Cells.Item(31, 13).Value = -1644.4736
$ws.Cells.Item(31, 14).Value = -6138.1665
# Row 34
$ws.Cells.Item(34, 8).Value = 2805.56
$ws.Cells.Item(34, 9).Value = 1939.4736
$ws.Cells.Item(34, 10).Value = 5548.1665
$ws.Cells.Item(34, 11).Value = 1939.4736
$ws.Cells.Item(34, 12).Value = 5548.1665
$ws.Cells.Item(34, 13).Value = -1737.4736
$ws.Cells.Item(34, 14).Value = -5952.1665
# Row 62
$ws.Cells.Item(62, 8).Value = 60514.855
$ws.Cells.Item(62, 10).Value = 103499.75
$ws.Cells.Item(62, 12).Value = 103499.75
$ws.Cells.Item(62, 14).Value = -104747.75
# Row 65
$ws.Cells.Item(65, 8).Value = 60514.855
$ws.Cells.Item(65, 10).Value = 103499.75
$ws.Cells.Item(65, 12).Value = 517498.75
$ws.Cells.Item(65, 14).Value = -523738.75
# Row 76
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 13).Value = $null
# Row 79
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 13).Value = $null
# Row 92
$ws.Cells.Item(92, 8).Value = 26375
$ws.Cells.Item(92, 10).Value = 26375
$ws.Cells.Item(92, 12).Value = 26375
$ws.Cells.Item(92, 14).Value = -31367
# Row 99
$ws.Cells.Item(99, 8).Value = 15400.368
$ws.Cells.Item(99, 9).Value = 12332.429
$ws.Cells.Item(99, 11).Value = 12332.429
$ws.Cells.Item(99, 13).Value = -10834.429
# Row 126
$ws.Cells.Item(126, 8).Value = 15400.368
$ws.Cells.Item(126, 9).Value = 12332.429
$ws.Cells.Item(126, 11).Value = 36997.287
$ws.Cells.Item(126, 13).Value = -34527.287

$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Cells.Item(103, 8).Value = 94.75
$ws.Cells.Item(103, 9).Value = 94.75
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 284.25
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = 594.75
$ws.Cells.Item(103, 14).Value = $null
# Row 107
$ws.Cells.Item(107, 8).Value = 1599.4
$ws.Cells.Item(107, 10).Value = 1777.25
$ws.Cells.Item(107, 12).Value = 5331.75
$ws.Cells.Item(107, 14).Value = -9171.75
# Row 113
$ws.Cells.Item(113, 8).Value = 1399.6666
$ws.Cells.Item(113, 9).Value = 899
$ws.Cells.Item(113, 11).Value = 2697
$ws.Cells.Item(113, 13).Value = -527
# Row 128
$ws.Cells.Item(128, 8).Value = 1457895.4
$ws.Cells.Item(128, 9).Value = 1457895.4
$ws.Cells.Item(128, 11).Value = 4373686.199999999
$ws.Cells.Item(128, 13).Value = -4368706.199999999
# Row 132
$ws.Cells.Item(132, 8).Value = 3502.3333
$ws.Cells.Item(132, 9).Value = 4030.3333
$ws.Cells.Item(132, 11).Value = 36272.9997
$ws.Cells.Item(132, 13).Value = -33742.9997

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 493.81818
$ws.Cells.Item(2, 9).Value = 254.125
$ws.Cells.Item(2, 10).Value = 1133
$ws.Cells.Item(2, 11).Value = 254.125
$ws.Cells.Item(2, 12).Value = 1133
$ws.Cells.Item(2, 13).Value = -141.125
$ws.Cells.Item(2, 14).Value = -1359
# Row 53
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 14).Value = $null
# Row 102
$ws.Cells.Item(102, 8).Value = 2774.1428
$ws.Cells.Item(102, 9).Value = 2027.1333
$ws.Cells.Item(102, 10).Value = 4641.6665
$ws.Cells.Item(102, 11).Value = 2027.1333
$ws.Cells.Item(102, 12).Value = 4641.6665
$ws.Cells.Item(102, 13).Value = -405.1333
$ws.Cells.Item(102, 14).Value = -7885.6665
# Row 122
$ws.Cells.Item(122, 8).Value = 86438.336
$ws.Cells.Item(122, 9).Value = 2184.7144
$ws.Cells.Item(122, 11).Value = 6554.1432
$ws.Cells.Item(122, 13).Value = -4104.1432
# Row 126
$ws.Cells.Item(126, 8).Value = 4502
$ws.Cells.Item(126, 10).Value = 4223.6665
$ws.Cells.Item(126, 12).Value = 12670.9995
$ws.Cells.Item(126, 14).Value = -17610.9995
# Row 132
$ws.Cells.Item(132, 8).Value = 1361.1428
$ws.Cells.Item(132, 9).Value = 840.1818
$ws.Cells.Item(132, 11).Value = 2520.5454
$ws.Cells.Item(132, 13).Value = 9.454600000000028

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 6961.154
$ws.Cells.Item(22, 9).Value = 3812.375
$ws.Cells.Item(22, 11).Value = 3812.375
$ws.Cells.Item(22, 13).Value = -3517.375
# Row 27
$ws.Cells.Item(27, 8).Value = 6961.154
$ws.Cells.Item(27, 9).Value = 3812.375
$ws.Cells.Item(27, 11).Value = 3812.375
$ws.Cells.Item(27, 13).Value = -3705.375
# Row 41
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).Value = $null
# Row 43
$ws.Cells.Item(43, 8).Value = 2907713.2
$ws.Cells.Item(43, 9).Value = 7998
$ws.Cells.Item(43, 10).Value = 5082500
$ws.Cells.Item(43, 11).Value = 7998
$ws.Cells.Item(43, 12).Value = 5082500
$ws.Cells.Item(43, 13).Value = -7805
$ws.Cells.Item(43, 14).Value = -5082886
# Row 46
$ws.Cells.Item(46, 8).Value = 5250
$ws.Cells.Item(46, 9).Value = 5000
$ws.Cells.Item(46, 10).Value = 5500
$ws.Cells.Item(46, 11).Value = 5000
$ws.Cells.Item(46, 12).Value = 5500
$ws.Cells.Item(46, 13).Value = -4812
$ws.Cells.Item(46, 14).Value = -5876
# Row 55
$ws.Cells.Item(55, 8).Value = 989.4545000000001
$ws.Cells.Item(55, 9).Value = 653
$ws.Cells.Item(55, 10).Value = 1393.2
$ws.Cells.Item(55, 11).Value = 653
$ws.Cells.Item(55, 12).Value = 1393.2
$ws.Cells.Item(55, 13).Value = -480
$ws.Cells.Item(55, 14).Value = -1739.2
# Row 122
$ws.Cells.Item(122, 8).Value = 8742
$ws.Cells.Item(122, 9).Value = 10489.333
$ws.Cells.Item(122, 10).Value = 3500
$ws.Cells.Item(122, 11).Value = 31467.999
$ws.Cells.Item(122, 12).Value = 10500
$ws.Cells.Item(122, 13).Value = -29017.999
$ws.Cells.Item(122, 14).Value = -15400
# Row 136
$ws.Cells.Item(136, 8).Value = 7058.375
$ws.Cells.Item(136, 9).Value = 6613.5
$ws.Cells.Item(136, 11).Value = 19840.5
$ws.Cells.Item(136, 13).Value = -17290.5

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Cells.Item(4, 8).Value = 2864641.5
$ws.Cells.Item(4, 9).Value = 6683496.5
$ws.Cells.Item(4, 10).Value = 500
$ws.Cells.Item(4, 11).Value = 6683496.5
$ws.Cells.Item(4, 12).Value = 500
$ws.Cells.Item(4, 13).Value = -6683383.5
$ws.Cells.Item(4, 14).Value = -726
# Row 126
$ws.Cells.Item(126, 8).Value = 2957.6155
$ws.Cells.Item(126, 9).Value = 2994.4443
$ws.Cells.Item(126, 11).Value = 8983.332900000001
$ws.Cells.Item(126, 13).Value = -6513.332900000001
# Row 132
$ws.Cells.Item(132, 8).Value = 2961.7576
$ws.Cells.Item(132, 9).Value = 2641.36
$ws.Cells.Item(132, 11).Value = 7924.08
$ws.Cells.Item(132, 13).Value = -5394.08
